# ---------------------------------------------------------------------------
# Adds a new "2022-Q4" sheet (right after "总计") with fund-holding data, and
# updates the "总计" (totals) sheet with a new row for 2022-Q4, shifting the
# existing 2022-Q3 / 2022-Q2 rows down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计" (2nd tab overall),
#    so the final tab order is: 总计, 2022-Q4, 2022-Q3, 2022-Q2.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q3SheetRef = $wb.Worksheets.Item("2022-Q3")

$q4Sheet = $wb.Worksheets.Add($q3SheetRef)
$q4Sheet.Name = "2022-Q4"

# NOTE: the worksheet object passed as the "Before" argument above becomes a
# stale reference once the new sheet has been inserted, so re-fetch a live
# handle on "2022-Q3" for everything that follows.
$q3Sheet = $wb.Worksheets.Item("2022-Q3")

# ---------------------------------------------------------------------------
# 2. Populate "2022-Q4" with the fund-holding table.
#    Columns B, D, E, F, G hold decimal-looking values that must stay exact
#    text (trailing zeros such as "0.6340" would be lost as numbers), so the
#    data rows (2-14) of those columns are pre-formatted as Text.
#    Columns A (row index) and H (rank) are plain numbers.
# ---------------------------------------------------------------------------
$q4Sheet.Range("B2:B14").NumberFormat = "@"
$q4Sheet.Range("D2:G14").NumberFormat = "@"

# Header row
$q4Sheet.Cells.Item(1,2).Value = "基金代码"
$q4Sheet.Cells.Item(1,3).Value = "基金名称"
$q4Sheet.Cells.Item(1,4).Value = "基金规模"
$q4Sheet.Cells.Item(1,5).Value = "股票总仓位"
$q4Sheet.Cells.Item(1,6).Value = "仓位占比"
$q4Sheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4Sheet.Cells.Item(1,8).Value = "仓位排名"

# row 2: 016067
$q4Sheet.Cells.Item(2,1).Value = 0
$q4Sheet.Cells.Item(2,2).Value = "016067"
$q4Sheet.Cells.Item(2,3).Value = "鹏华新能源汽车混合A"
$q4Sheet.Cells.Item(2,4).Value = "19.04"
$q4Sheet.Cells.Item(2,5).Value = "95.25"
$q4Sheet.Cells.Item(2,6).Value = "3.33"
$q4Sheet.Cells.Item(2,7).Value = "0.6340"
$q4Sheet.Cells.Item(2,8).Value = 8
# row 3: 160211
$q4Sheet.Cells.Item(3,1).Value = 1
$q4Sheet.Cells.Item(3,2).Value = "160211"
$q4Sheet.Cells.Item(3,3).Value = "国泰中小盘成长混合（LOF）"
$q4Sheet.Cells.Item(3,4).Value = "6.59"
$q4Sheet.Cells.Item(3,5).Value = "90.98"
$q4Sheet.Cells.Item(3,6).Value = "2.16"
$q4Sheet.Cells.Item(3,7).Value = "0.1423"
$q4Sheet.Cells.Item(3,8).Value = 6
# row 4: 014686
$q4Sheet.Cells.Item(4,1).Value = 2
$q4Sheet.Cells.Item(4,2).Value = "014686"
$q4Sheet.Cells.Item(4,3).Value = "招商核心装备混合A"
$q4Sheet.Cells.Item(4,4).Value = "1.90"
$q4Sheet.Cells.Item(4,5).Value = "91.32"
$q4Sheet.Cells.Item(4,6).Value = "6.52"
$q4Sheet.Cells.Item(4,7).Value = "0.1239"
$q4Sheet.Cells.Item(4,8).Value = 1
# row 5: 016068
$q4Sheet.Cells.Item(5,1).Value = 3
$q4Sheet.Cells.Item(5,2).Value = "016068"
$q4Sheet.Cells.Item(5,3).Value = "鹏华新能源汽车混合C"
$q4Sheet.Cells.Item(5,4).Value = "3.65"
$q4Sheet.Cells.Item(5,5).Value = "95.25"
$q4Sheet.Cells.Item(5,6).Value = "3.33"
$q4Sheet.Cells.Item(5,7).Value = "0.1215"
$q4Sheet.Cells.Item(5,8).Value = 8
# row 6: 014320
$q4Sheet.Cells.Item(6,1).Value = 4
$q4Sheet.Cells.Item(6,2).Value = "014320"
$q4Sheet.Cells.Item(6,3).Value = "德邦半导体产业混合C"
$q4Sheet.Cells.Item(6,4).Value = "1.52"
$q4Sheet.Cells.Item(6,5).Value = "92.57"
$q4Sheet.Cells.Item(6,6).Value = "4.60"
$q4Sheet.Cells.Item(6,7).Value = "0.0699"
$q4Sheet.Cells.Item(6,8).Value = 8
# row 7: 014687
$q4Sheet.Cells.Item(7,1).Value = 5
$q4Sheet.Cells.Item(7,2).Value = "014687"
$q4Sheet.Cells.Item(7,3).Value = "招商核心装备混合C"
$q4Sheet.Cells.Item(7,4).Value = "0.84"
$q4Sheet.Cells.Item(7,5).Value = "91.32"
$q4Sheet.Cells.Item(7,6).Value = "6.52"
$q4Sheet.Cells.Item(7,7).Value = "0.0548"
$q4Sheet.Cells.Item(7,8).Value = 1
# row 8: 001815
$q4Sheet.Cells.Item(8,1).Value = 6
$q4Sheet.Cells.Item(8,2).Value = "001815"
$q4Sheet.Cells.Item(8,3).Value = "华泰柏瑞激励动力灵活配置混合A"
$q4Sheet.Cells.Item(8,4).Value = "1.96"
$q4Sheet.Cells.Item(8,5).Value = "88.25"
$q4Sheet.Cells.Item(8,6).Value = "2.24"
$q4Sheet.Cells.Item(8,7).Value = "0.0439"
$q4Sheet.Cells.Item(8,8).Value = 8
# row 9: 002082
$q4Sheet.Cells.Item(9,1).Value = 7
$q4Sheet.Cells.Item(9,2).Value = "002082"
$q4Sheet.Cells.Item(9,3).Value = "华泰柏瑞激励动力灵活配置混合C"
$q4Sheet.Cells.Item(9,4).Value = "0.95"
$q4Sheet.Cells.Item(9,5).Value = "88.25"
$q4Sheet.Cells.Item(9,6).Value = "2.24"
$q4Sheet.Cells.Item(9,7).Value = "0.0213"
$q4Sheet.Cells.Item(9,8).Value = 8
# row 10: 010571
$q4Sheet.Cells.Item(10,1).Value = 8
$q4Sheet.Cells.Item(10,2).Value = "010571"
$q4Sheet.Cells.Item(10,3).Value = "新沃创新领航混合C"
$q4Sheet.Cells.Item(10,4).Value = "0.51"
$q4Sheet.Cells.Item(10,5).Value = "93.56"
$q4Sheet.Cells.Item(10,6).Value = "3.71"
$q4Sheet.Cells.Item(10,7).Value = "0.0189"
$q4Sheet.Cells.Item(10,8).Value = 10
# row 11: 014319
$q4Sheet.Cells.Item(11,1).Value = 9
$q4Sheet.Cells.Item(11,2).Value = "014319"
$q4Sheet.Cells.Item(11,3).Value = "德邦半导体产业混合A"
$q4Sheet.Cells.Item(11,4).Value = "0.37"
$q4Sheet.Cells.Item(11,5).Value = "92.57"
$q4Sheet.Cells.Item(11,6).Value = "4.60"
$q4Sheet.Cells.Item(11,7).Value = "0.0170"
$q4Sheet.Cells.Item(11,8).Value = 8
# row 12: 010570
$q4Sheet.Cells.Item(12,1).Value = 10
$q4Sheet.Cells.Item(12,2).Value = "010570"
$q4Sheet.Cells.Item(12,3).Value = "新沃创新领航混合A"
$q4Sheet.Cells.Item(12,4).Value = "0.24"
$q4Sheet.Cells.Item(12,5).Value = "93.56"
$q4Sheet.Cells.Item(12,6).Value = "3.71"
$q4Sheet.Cells.Item(12,7).Value = "0.0089"
$q4Sheet.Cells.Item(12,8).Value = 10
# row 13: 012143
$q4Sheet.Cells.Item(13,1).Value = 11
$q4Sheet.Cells.Item(13,2).Value = "012143"
$q4Sheet.Cells.Item(13,3).Value = "新沃内需增长混合A"
$q4Sheet.Cells.Item(13,4).Value = "0.20"
$q4Sheet.Cells.Item(13,5).Value = "93.63"
$q4Sheet.Cells.Item(13,6).Value = "3.96"
$q4Sheet.Cells.Item(13,7).Value = "0.0079"
$q4Sheet.Cells.Item(13,8).Value = 10
# row 14: 012144
$q4Sheet.Cells.Item(14,1).Value = 12
$q4Sheet.Cells.Item(14,2).Value = "012144"
$q4Sheet.Cells.Item(14,3).Value = "新沃内需增长混合C"
$q4Sheet.Cells.Item(14,4).Value = "0.04"
$q4Sheet.Cells.Item(14,5).Value = "93.63"
$q4Sheet.Cells.Item(14,6).Value = "3.96"
$q4Sheet.Cells.Item(14,7).Value = "0.0016"
$q4Sheet.Cells.Item(14,8).Value = 10

# Match the bold/bordered look of the header row and the row-index column,
# copying the format straight from the sibling "2022-Q3" sheet.
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q3Sheet.Range("A2").Copy()
$q4Sheet.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Update the "总计" sheet: shift 2022-Q3 / 2022-Q2 rows down and insert the
#    new 2022-Q4 totals in row 2.
# ---------------------------------------------------------------------------
$totalSheet.Cells.Item(3,1).Copy()
$totalSheet.Cells.Item(4,1).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$totalSheet.Cells.Item(4,1).Value = 2
$totalSheet.Cells.Item(4,2).Value = "2022-Q2"
$totalSheet.Cells.Item(4,3).Value = 14
$totalSheet.Cells.Item(4,4).Value = 3.45

$totalSheet.Cells.Item(3,1).Value = 1
$totalSheet.Cells.Item(3,2).Value = "2022-Q3"
$totalSheet.Cells.Item(3,3).Value = 15
$totalSheet.Cells.Item(3,4).Value = 0.65

$totalSheet.Cells.Item(2,1).Value = 0
$totalSheet.Cells.Item(2,2).Value = "2022-Q4"
$totalSheet.Cells.Item(2,3).Value = 13
$totalSheet.Cells.Item(2,4).Value = 1.27

# ---------------------------------------------------------------------------
# 4. Restore "2022-Q2" as the selected/active tab (adding the new sheet made
#    it active by default).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Select()
